$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 8 formatting down to the two new rows (9 and 10) before writing values,
# so the bold/border style on column A carries over to the newly appended rows.
$ws.Range("A8").Copy()
$ws.Range("A9:A10").PasteSpecial(-4122)

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2023
$ws.Range("D2").Value = 237
$ws.Range("F2").Value = 112.8173076923077
$ws.Range("G2").Value = 6.5
$ws.Range("H2").Value = "Detroit"
$ws.Range("I2").Value = "Washington"
$ws.Range("J2").Value = 0.5238970588235294
$ws.Range("K2").Value = 99.06357692307694
$ws.Range("L2").Value = 113.4078846153846
$ws.Range("M2").Value = 117.0733076923077
$ws.Range("N2").Value = 75.60203846153846
$ws.Range("O2").Value = 0.371936923076923
$ws.Range("P2").Value = 0.5733776923076923
$ws.Range("Q2").Value = 0.2925934615384615
$ws.Range("R2").Value = 12.24257692307692
$ws.Range("S2").Value = 11.311
$ws.Range("T2").Value = 0.2222367307692307
$ws.Range("U2").Value = 0.9878923615788765
$ws.Range("V2").Value = 1.029693435936023
$ws.Range("W2").Value = 11.05513080580208
$ws.Range("X2").Value = 0.365
$ws.Range("Y2").Value = 32.5
$ws.Range("Z2").Value = 75.15
$ws.Range("AA2").Value = 0.4980185443817948

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2023
$ws.Range("D3").Value = 239
$ws.Range("F3").Value = 115.39
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = "Memphis"
$ws.Range("I3").Value = "Portland"
$ws.Range("J3").Value = 0.4583333333333333
$ws.Range("K3").Value = 98.94099999999997
$ws.Range("L3").Value = 116.205
$ws.Range("M3").Value = 113.935
$ws.Range("N3").Value = 76.00299999999999
$ws.Range("O3").Value = 0.38024
$ws.Range("P3").Value = 0.57931
$ws.Range("Q3").Value = 0.28796
$ws.Range("R3").Value = 12.209
$ws.Range("S3").Value = 12.132
$ws.Range("T3").Value = 0.214475
$ws.Range("U3").Value = 1.010420315236427
$ws.Range("V3").Value = 1.011758781188778
$ws.Range("W3").Value = 12.50804112607605
$ws.Range("X3").Value = 0.56
$ws.Range("Y3").Value = 44.5
$ws.Range("Z3").Value = 75.15
$ws.Range("AA3").Value = 0.520528442858359

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2023
$ws.Range("D4").Value = 230.5
$ws.Range("F4").Value = 112.9641856742697
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = "Philadelphia"
$ws.Range("I4").Value = "Orlando"
$ws.Range("J4").Value = 0.5714285714285714
$ws.Range("K4").Value = 97.77575030012007
$ws.Range("L4").Value = 114.7935774309724
$ws.Range("M4").Value = 114.627130852341
$ws.Range("N4").Value = 77.0796318527411
$ws.Range("O4").Value = 0.3807599039615847
$ws.Range("P4").Value = 0.591280712284914
$ws.Range("Q4").Value = 0.2958035214085634
$ws.Range("R4").Value = 12.70594237695078
$ws.Range("S4").Value = 12.57018807523009
$ws.Range("T4").Value = 0.2273242296918767
$ws.Range("U4").Value = 0.9891785085312584
$ws.Range("V4").Value = 1.039520022548001
$ws.Range("W4").Value = 10.89118581312156
$ws.Range("X4").Value = 0.522609043617447
$ws.Range("Y4").Value = 38.5
$ws.Range("Z4").Value = 75.4
$ws.Range("AA4").Value = 0.4851847805777361

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 2023
$ws.Range("D5").Value = 226.5
$ws.Range("F5").Value = 116.1623529411765
$ws.Range("G5").Value = 8.5
$ws.Range("H5").Value = "Boston"
$ws.Range("I5").Value = "Brooklyn"
$ws.Range("J5").Value = 0.4591836734693877
$ws.Range("K5").Value = 97.98670588235294
$ws.Range("L5").Value = 117.7425294117647
$ws.Range("M5").Value = 113.3143333333333
$ws.Range("N5").Value = 76.20352941176469
$ws.Range("O5").Value = 0.4242641176470587
$ws.Range("P5").Value = 0.6076550980392157
$ws.Range("Q5").Value = 0.2537174509803922
$ws.Range("R5").Value = 12.23876470588236
$ws.Range("S5").Value = 11.34262745098039
$ws.Range("T5").Value = 0.2031652941176471
$ws.Range("U5").Value = 1.017183475842176
$ws.Range("V5").Value = 1.007956952265612
$ws.Range("W5").Value = 11.22529722880983
$ws.Range("X5").Value = 0.6629411764705883
$ws.Range("Y5").Value = 52.5
$ws.Range("Z5").Value = 77.55
$ws.Range("AA5").Value = 0.4922572198882214

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 2023
$ws.Range("D6").Value = 238
$ws.Range("F6").Value = 116.5877358490566
$ws.Range("G6").Value = 3.5
$ws.Range("H6").Value = "Minnesota"
$ws.Range("I6").Value = "GoldenState"
$ws.Range("J6").Value = 0.528498427672956
$ws.Range("K6").Value = 101.1109433962264
$ws.Range("L6").Value = 114.6936981132076
$ws.Range("M6").Value = 114.6116603773585
$ws.Range("N6").Value = 74.86296226415094
$ws.Range("O6").Value = 0.4322760377358491
$ws.Range("P6").Value = 0.5970496226415093
$ws.Range("Q6").Value = 0.2525354716981132
$ws.Range("R6").Value = 13.7728679245283
$ws.Range("S6").Value = 12.70501886792453
$ws.Range("T6").Value = 0.21137
$ws.Range("U6").Value = 1.020908369956713
$ws.Range("V6").Value = 1.025839593642522
$ws.Range("W6").Value = 10.71537865251756
$ws.Range("X6").Value = 0.5147169811320755
$ws.Range("Y6").Value = 51
$ws.Range("Z6").Value = 76
$ws.Range("AA6").Value = 0.4741815263775075

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 2023
$ws.Range("D7").Value = 232.5
$ws.Range("F7").Value = 113.22
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = "Houston"
$ws.Range("I7").Value = "OklahomaCity"
$ws.Range("J7").Value = 0.53
$ws.Range("K7").Value = 99.66099999999997
$ws.Range("L7").Value = 112.65
$ws.Range("M7").Value = 115.993
$ws.Range("N7").Value = 74.42599999999999
$ws.Range("O7").Value = 0.37297
$ws.Range("P7").Value = 0.5612999999999999
$ws.Range("Q7").Value = 0.26752
$ws.Range("R7").Value = 12.744
$ws.Range("S7").Value = 12.703
$ws.Range("T7").Value = 0.214635
$ws.Range("U7").Value = 0.9914185639229421
$ws.Range("V7").Value = 0.9982417372526393
$ws.Range("W7").Value = 11.35944194008727
$ws.Range("X7").Value = 0.36
$ws.Range("Y7").Value = 23.5
$ws.Range("Z7").Value = 74.05000000000001
$ws.Range("AA7").Value = 0.5114597516781482

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 2023
$ws.Range("D8").Value = 244
$ws.Range("F8").Value = 116.0492196878752
$ws.Range("G8").Value = 7.5
$ws.Range("H8").Value = "SanAntonio"
$ws.Range("I8").Value = "Sacramento"
$ws.Range("J8").Value = 0.5352891156462585
$ws.Range("K8").Value = 100.4937775110044
$ws.Range("L8").Value = 115.121168467387
$ws.Range("M8").Value = 118.5173269307723
$ws.Range("N8").Value = 76.5344137655062
$ws.Range("O8").Value = 0.3792406962785113
$ws.Range("P8").Value = 0.5837677070828331
$ws.Range("Q8").Value = 0.2631828731492597
$ws.Range("R8").Value = 12.58831532613045
$ws.Range("S8").Value = 12.01444577831132
$ws.Range("T8").Value = 0.2011364545818328
$ws.Range("U8").Value = 1.016192816881569
$ws.Range("V8").Value = 0.9300735056900263
$ws.Range("W8").Value = 11.37890465049901
$ws.Range("X8").Value = 0.42296918767507
$ws.Range("Y8").Value = 28.5
$ws.Range("Z8").Value = 74.55000000000001
$ws.Range("AA8").Value = 0.470376252904394

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 2023
$ws.Range("D9").Value = 230.5
$ws.Range("F9").Value = 115.1442307692308
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = "Utah"
$ws.Range("I9").Value = "Toronto"
$ws.Range("J9").Value = 0.5533559577677225
$ws.Range("K9").Value = 97.84134615384616
$ws.Range("L9").Value = 116.9346153846154
$ws.Range("M9").Value = 116.1682692307692
$ws.Range("N9").Value = 75.91346153846153
$ws.Range("O9").Value = 0.4012596153846154
$ws.Range("P9").Value = 0.5740480769230769
$ws.Range("Q9").Value = 0.2747788461538461
$ws.Range("R9").Value = 11.45
$ws.Range("S9").Value = 13.09230769230769
$ws.Range("T9").Value = 0.2177884615384615
$ws.Range("U9").Value = 1.008268220396066
$ws.Range("V9").Value = 1.010572394236757
$ws.Range("W9").Value = 10.32440218516663
$ws.Range("X9").Value = 0.4711538461538461
$ws.Range("Y9").Value = 35
$ws.Range("Z9").Value = 74.6
$ws.Range("AA9").Value = 0.4974347709542146

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 2023
$ws.Range("D10").Value = 229.5
$ws.Range("F10").Value = 114.5497737556561
$ws.Range("G10").Value = 1.5
$ws.Range("H10").Value = "Phoenix"
$ws.Range("I10").Value = "Atlanta"
$ws.Range("J10").Value = 0.5392156862745098
$ws.Range("K10").Value = 98.45803167420814
$ws.Range("L10").Value = 115.3364819004525
$ws.Range("M10").Value = 114.9615007541478
$ws.Range("N10").Value = 75.66872171945701
$ws.Range("O10").Value = 0.3516091628959275
$ws.Range("P10").Value = 0.57205580693816
$ws.Range("Q10").Value = 0.2406468702865762
$ws.Range("R10").Value = 11.39298642533937
$ws.Range("S10").Value = 12.3196455505279
$ws.Range("T10").Value = 0.2091603506787331
$ws.Range("U10").Value = 1.003062817475097
$ws.Range("V10").Value = 1.035363421055657
$ws.Range("W10").Value = 11.81077167049513
$ws.Range("X10").Value = 0.5047134238310709
$ws.Range("Y10").Value = 49.5
$ws.Range("Z10").Value = 76.9
$ws.Range("AA10").Value = 0.500478043140146

